$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = "2026-02-16 17:48:59"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "97%"
$ws.Range("I2").Value = "18.0 mm"
$ws.Range("E3").Value = "2026-02-16 17:49:02"
$ws.Range("G3").Value = "231 cm"
$ws.Range("I3").Value = "9.5 mm"
$ws.Range("L3").Value = "60.1 km/h - 260º 17:17 TU"
$ws.Range("N3").Value = "-2.1 °C 17:04 TU"
$ws.Range("E4").Value = "2026-02-16 17:49:05"
$ws.Range("E5").Value = "2026-02-16 17:49:08"
$ws.Range("G5").Value = "145 cm"
$ws.Range("I5").Value = "21.5 mm"
$ws.Range("K5").Value = "4.1 MJ/m2"
$ws.Range("L5").Value = "41.0 km/h - 321º 17:26 TU"
$ws.Range("N5").Value = "-1.5 °C 17:20 TU"
$ws.Range("E6").Value = "2026-02-16 17:49:10"
$ws.Range("O6").Value = "11.6 °C"
$ws.Range("E7").Value = "2026-02-16 17:49:13"
$ws.Range("O7").Value = "16.2 °C"
$ws.Range("E8").Value = "2026-02-16 17:49:16"
$ws.Range("K8").Value = "12.1 MJ/m2"
$ws.Range("O8").Value = "12.4 °C"
$ws.Range("E9").Value = "2026-02-16 17:49:19"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "73%"
$ws.Range("O9").Value = "11.0 °C"
$ws.Range("E10").Value = "2026-02-16 17:49:22"
$ws.Range("O10").Value = "10.8 °C"
$ws.Range("E11").Value = "2026-02-16 17:49:24"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "79%"
$ws.Range("O11").Value = "6.5 °C"
$ws.Range("E12").Value = "2026-02-16 17:49:27"
$ws.Range("O12").Value = "10.4 °C"
$ws.Range("E13").Value = "2026-02-16 17:49:30"
$ws.Range("J13").Value = "1015.0 hPa"
$ws.Range("O13").Value = "5.5 °C"
$ws.Range("E14").Value = "2026-02-16 17:49:33"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "56%"
$ws.Range("O14").Value = "15.9 °C"
$ws.Range("E15").Value = "2026-02-16 17:49:35"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "70%"
$ws.Range("O15").Value = "11.1 °C"
$ws.Range("E16").Value = "2026-02-16 17:49:38"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "79%"
$ws.Range("E17").Value = "2026-02-16 17:49:41"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "68%"
$ws.Range("N17").Value = "4.3 °C 17:29 TU"
$ws.Range("O17").Value = "6.2 °C"
$ws.Range("E18").Value = "2026-02-16 17:49:44"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "74%"
$ws.Range("J18").Value = "1012.9 hPa"
$ws.Range("O18").Value = "10.8 °C"
$ws.Range("E19").Value = "2026-02-16 17:49:46"
$ws.Range("O19").Value = "6.8 °C"
$ws.Range("E20").Value = "2026-02-16 17:49:49"
$ws.Range("I20").Value = "0.1 mm"
$ws.Range("E21").Value = "2026-02-16 17:49:52"
$ws.Range("J21").Value = "1014.5 hPa"
$ws.Range("O21").Value = "8.2 °C"
$ws.Range("E22").Value = "2026-02-16 17:49:55"
$ws.Range("E23").Value = "2026-02-16 17:49:58"
$ws.Range("I23").Value = "12.4 mm"
$ws.Range("N23").Value = "-1.9 °C 17:19 TU"
$ws.Range("E24").Value = "2026-02-16 17:50:00"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "70%"
$ws.Range("J24").Value = "1016.8 hPa"
$ws.Range("E25").Value = "2026-02-16 17:50:03"
$ws.Range("I25").Value = "5.2 mm"
$ws.Range("E26").Value = "2026-02-16 17:50:06"
$ws.Range("E27").Value = "2026-02-16 17:50:08"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "82%"
$ws.Range("L27").Value = "39.6 km/h - 216º 17:21 TU"
$ws.Range("O27").Value = "1.3 °C"
$ws.Range("E28").Value = "2026-02-16 17:50:11"
$ws.Range("J28").Value = "1013.0 hPa"
$ws.Range("O28").Value = "9.4 °C"
$ws.Range("E29").Value = "2026-02-16 17:50:14"
$ws.Range("E30").Value = "2026-02-16 17:50:17"
$ws.Range("J30").Value = "1012.5 hPa"
$ws.Range("E31").Value = "2026-02-16 17:50:19"
$ws.Range("J31").Value = "1011.7 hPa"
$ws.Range("E32").Value = "2026-02-16 17:50:22"
$ws.Range("E33").Value = "2026-02-16 17:50:25"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "73%"
$ws.Range("E34").Value = "2026-02-16 17:50:27"
$ws.Range("L34").Value = "51.8 km/h - 31º 17:05 TU"
$ws.Range("E35").Value = "2026-02-16 17:50:30"
$ws.Range("K35").Value = "11.8 MJ/m2"
$ws.Range("E36").Value = "2026-02-16 17:50:33"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "75%"
$ws.Range("J36").Value = "1012.8 hPa"
$ws.Range("O36").Value = "11.5 °C"
$ws.Range("E37").Value = "2026-02-16 17:50:36"
$ws.Range("J37").Value = "1014.9 hPa"
$ws.Range("O37").Value = "6.4 °C"
$ws.Range("E38").Value = "2026-02-16 17:50:39"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "71%"
$ws.Range("O38").Value = "11.8 °C"
$ws.Range("E39").Value = "2026-02-16 17:50:42"
$ws.Range("E40").Value = "2026-02-16 17:50:44"
$ws.Range("J40").Value = "1016.7 hPa"
$ws.Range("O40").Value = "6.7 °C"
$ws.Range("E41").Value = "2026-02-16 17:50:47"
$ws.Range("J41").Value = "1015.0 hPa"
$ws.Range("O41").Value = "17.4 °C"
$ws.Range("E42").Value = "2026-02-16 17:50:50"
$ws.Range("O42").Value = "11.1 °C"
$ws.Range("E43").Value = "2026-02-16 17:50:52"
$ws.Range("O43").Value = "8.1 °C"
$ws.Range("E44").Value = "2026-02-16 17:50:55"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "89%"
$ws.Range("K44").Value = "8.5 MJ/m2"
$ws.Range("N44").Value = "-1.8 °C 17:29 TU"
$ws.Range("E45").Value = "2026-02-16 17:50:58"
$ws.Range("I45").Value = "14.9 mm"
$ws.Range("E46").Value = "2026-02-16 17:51:00"
$ws.Range("J46").Value = "1017.1 hPa"
$ws.Range("K46").Value = "12.8 MJ/m2"
$ws.Range("O46").Value = "15.9 °C"
